$wb = $excel.ActiveWorkbook

# --- Add the new "CreateJobTestData" worksheet after the existing one ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CreateJobTestData"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mst_service_location_id"
$ws.Range("B1").Value = "mst_platform_id"
$ws.Range("C1").Value = "mst_warrenty_status_id"
$ws.Range("D1").Value = "mst_oem_id"
$ws.Range("E1").Value = "customer__first_name"
$ws.Range("F1").Value = "customer__last_name"
$ws.Range("G1").Value = "customer__mobile_number"
$ws.Range("H1").Value = "customer__mobile_number_alt"
$ws.Range("I1").Value = "customer__email_id"
$ws.Range("J1").Value = "customer__email_id_alt"
$ws.Range("K1").Value = "customer_address__flat_number"
$ws.Range("L1").Value = "customer_address__apartment_name"
$ws.Range("M1").Value = "customer_address__street_name"
$ws.Range("N1").Value = "customer_address__landmark"
$ws.Range("O1").Value = "customer_address__area"
$ws.Range("P1").Value = "customer_address__pincode"
$ws.Range("Q1").Value = "customer_address__country"
$ws.Range("R1").Value = "customer_address__state"
$ws.Range("S1").Value = "customer_product__dop"
$ws.Range("T1").Value = "customer_product__serial_number"
$ws.Range("U1").Value = "customer_product__imei1"
$ws.Range("V1").Value = "customer_product__imei2"
$ws.Range("W1").Value = "customer_product__popurl"
$ws.Range("X1").Value = "customer_product__product_id"
$ws.Range("Y1").Value = "customer_product__mst_model_id"
$ws.Range("Z1").Value = "problems__id"
$ws.Range("AA1").Value = "problems__remark"

# --- Data row (row 2), left-to-right, except the serial/imei/alt-email
#     cells which are keyed in afterwards (see below) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Hilda"
$ws.Range("F2").Value = "Sipes"
$ws.Range("G2").Value = 5125173023
$ws.Range("H2").Value = "512-517-3023"
$ws.Range("I2").Value = "Nakia87@hotmail.com"
$ws.Range("J2").Value = "Nakia87@hotmail.com"
$ws.Range("K2").Value = 768
$ws.Range("L2").Value = "Swamy nagar"
$ws.Range("M2").Value = "61638 Graham Passage"
$ws.Range("N2").Value = "Near Shivalayam"
$ws.Range("O2").Value = "Kakinada"
$ws.Range("P2").Value = 516
$ws.Range("Q2").Value = "Faroe Islands"
$ws.Range("R2").Value = "Andhra Pradesh"
$ws.Range("S2").Value = "2025-09-30T18:30:00.000Z"
$ws.Range("W2").Value = "2025-09-30T18:30:00.000Z"
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 6
$ws.Range("AA2").Value = "over heat"

# customer_product serial/imei numbers: keyed in last, forced to text with a
# leading apostrophe (quote-prefix) and shown in scientific notation -
# matches how these long numeric-looking identifiers get entered in Excel.
$ws.Range("T2").Value = "'12345700000000"
$ws.Range("T2").NumberFormat = "0.00E+00"
$ws.Range("U2").Value = "'96257600000000"
$ws.Range("U2").NumberFormat = "0.00E+00"
$ws.Range("V2").Value = "'96257600000000"
$ws.Range("V2").NumberFormat = "0.00E+00"

# customer__email_id_alt re-keyed last as forced text (quote-prefix)
$ws.Range("J2").Value = "'Nakia87@hotmail.com"

# --- Column widths tweaked on this sheet ---
$ws.Columns.Item(9).ColumnWidth = 12.833333333333334
$ws.Columns.Item(10).ColumnWidth = 15.333333333333334
$ws.Columns.Item(19).ColumnWidth = 13.666666666666666

# --- Selection on the new sheet, which becomes the active tab ---
$ws.Range("K10").Select() | Out-Null
